$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns with refreshed crypto data.
# D-column values are forced to Text so Excel does not reinterpret
# dotted-thousands price strings (e.g. "37.640.93") as numbers/dates,
# then the cell style is reset to "Normal" so no stray number format
# is left behind on cells that started with the default style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.640.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.067.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0779"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.375.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.767"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.088.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.586.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.137"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.72%  "

$ws.Range("E36").Value = "  +2.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0971"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  -2.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.452.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.259.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
